$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.517.54"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").Value = "2.633.20"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.648"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.61%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -3.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.82"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("E11").Value = "  -2.11%  "

$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("E14").Value = "  -5.40%  "

$ws.Range("D15").Value = "3.106.03"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").Value = "64.323.12"
$ws.Range("E16").Value = "  -1.82%  "

$ws.Range("D17").Value = "2.649.74"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.28"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.90%  "

$ws.Range("E19").Value = "  -2.31%  "

$ws.Range("E20").Value = "  -0.89%  "

$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.83%  "

$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("E25").Value = "  +4.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.31%  "

$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "566.20"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.64%  "

$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("E32").Value = "  -2.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.66"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.15%  "

$ws.Range("E34").Value = "  -3.53%  "

$ws.Range("E35").Value = "  -2.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.414"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.08"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.95"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.46"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "158.50"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.52%  "

$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.14"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.32%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0601"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.637"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("E48").Value = "  +3.98%  "

$ws.Range("E49").Value = "  -2.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.20"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.13%  "

$ws.Range("D51").Value = "0.0₆0240"
$ws.Range("E51").Value = "  -5.14%  "
